$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added at the top of the "Haba" price history
# (Femacal de La Calera). Insert a new row 69 and push the existing
# rows 69:99 down to 70:100, then populate the new row with the
# latest observation.
$ws.Rows.Item(69).Insert()

$ws.Cells.Item(69, 1).Value = 3
$ws.Cells.Item(69, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(69, 3).Value = "Coquimbo"
$ws.Cells.Item(69, 4).Value = 44523
$ws.Cells.Item(69, 4).NumberFormat = $ws.Cells.Item(70, 4).NumberFormat
$ws.Cells.Item(69, 5).Value = 5
$ws.Cells.Item(69, 6).Value = 100112026
$ws.Cells.Item(69, 7).Value = "Haba"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 110
$ws.Cells.Item(69, 11).Value = 7000
$ws.Cells.Item(69, 12).Value = 7500
$ws.Cells.Item(69, 13).Value = 7227
$ws.Cells.Item(69, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(69, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(69, 16).Value = 289
$ws.Cells.Item(69, 17).Value = 25
$ws.Cells.Item(69, 18).Value = "Hortaliza"
